$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Courses")

# Insert a new column before column I (9th column). Existing columns
# I, J, K shift right to become J, K, L. The new column picks up the
# formatting (style) of column H automatically.
$ws.Columns("I:I").Insert()

# Set the new header cell (I1) to the new "ANTI-REQUISITE" column title.
$ws.Range("I1").Value = "ANTI-REQUISITE"

# Copy the header formatting (fill/border/font/alignment) from the
# neighbouring header cell (J1, the old I1) onto I1 and H1 so that both
# match the centered header style used by the other header cells.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill the new ANTI-REQUISITE column with the default value "None" for
# every course data row (rows 2 through 22).
$ws.Range("I2:I22").Value = "None"

# Update the selection to reflect where the editor left off.
$ws.Range("J6").Select()
